$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1. Insert a new "Meta description" paragraph right after the first
#    (Heading1 title) paragraph.
#
#    To avoid inheriting the Heading1 paragraph style (and to avoid
#    the rsid stamping that an explicit Style re-assignment causes),
#    build the new paragraph by splitting the third paragraph - which
#    already carries the default/"Normal" (style-less) formatting -
#    right at its very start. That produces a clean, style-less
#    paragraph that we fill in and then relocate into place.
# ------------------------------------------------------------------
$splitTarget = $d.Paragraphs.Item(3)
$splitPos = $splitTarget.Range.Start
$rSplit = $d.Range($splitPos, $splitPos)
$rSplit.InsertParagraphBefore()

$cleanPara = $d.Paragraphs.Item(3)
$fillStart = $cleanPara.Range.Start

$boldRun = $d.Range($fillStart, $fillStart)
$boldRun.InsertAfter("Meta description")

$afterBold = $boldRun.End
$restRun = $d.Range($afterBold, $afterBold)
$restRun.InsertAfter(": Read our Blirix Workshop review and play for free! Enjoy the impressive Steampunk setting, high volatility, and advanced spells during free spins.")

# Apply bold formatting only to the "Meta description" label.
$boldRange = $d.Range($fillStart, $afterBold)
$boldRange.Font.Bold = 1

# Relocate the filled-in paragraph so it sits right after the title
# (Heading1) paragraph instead of right before "Symbols and Wilds".
$filledPara = $d.Paragraphs.Item(3)
$moveRange = $d.Range($filledPara.Range.Start, $filledPara.Range.End)
$moveRange.Cut()

$titlePara = $d.Paragraphs.Item(1)
$pastePos = $titlePara.Range.End
$pasteRange = $d.Range($pastePos, $pastePos)
$pasteRange.Paste()

# ------------------------------------------------------------------
# 2. Remove the old duplicate bold title paragraph further down in
#    the document (the one that is not the Heading1 paragraph).
# ------------------------------------------------------------------
$dupIdx = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    $t = $p.Range.Text.TrimEnd()
    if ($t -eq "Play Blirix Workshop Free - Impressive Steampunk Setting" -and $i -ne 1) {
        $dupIdx = $i
    }
}

if ($dupIdx -ne -1) {
    $dupPara = $d.Paragraphs.Item($dupIdx)
    $dupRange = $d.Range($dupPara.Range.Start, $dupPara.Range.End)
    $dupRange.Delete()
}

# ------------------------------------------------------------------
# 3. Replace the text of the final (italic) paragraph with the new
#    image-prompt copy, preserving its italic run formatting. Using
#    a direct Range.Text assignment (rather than Find/Replace) keeps
#    straight quotes/apostrophes from being "smart-quoted".
# ------------------------------------------------------------------
$lastIdx = $d.Paragraphs.Count
$lastPara = $d.Paragraphs.Item($lastIdx)
$lastRange = $d.Range($lastPara.Range.Start, $lastPara.Range.End)

$newClosingText = 'Please create a feature image for the game "Blirix Workshop" that features a happy Maya warrior with glasses. The image should be in cartoon style, and should be eye-catching and engaging. The Maya warrior should be smiling and holding a bubbling test tube, with the Blirix Workshop logo and some of the game''s symbols (such as the carnivorous plant, book of potions, and hourglass) featured around them. The background should also be steampunk-inspired, with gears and machinery visible. The image should convey a sense of excitement and adventure, and make players eager to explore the fantastical world of Blirix Workshop.'

$lastRange.Text = $newClosingText

Write-Output "done"
